$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 7841
$wsExhibit.Range("F15").Value = 9312
$wsExhibit.Range("F24").Value = 1213
$wsExhibit.Range("F25").Value = 126
$wsExhibit.Range("F27").Value = 745
$wsExhibit.Range("F28").Value = 966
$wsExhibit.Range("F30").Value = 1904
$wsExhibit.Range("F31").Value = 348
$wsExhibit.Range("F34").Value = 1499

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 9312
$wsAll.Range("F23").Value = 1213
$wsAll.Range("F24").Value = 126
$wsAll.Range("F26").Value = 745
$wsAll.Range("F27").Value = 966
$wsAll.Range("F29").Value = 1904
$wsAll.Range("F30").Value = 348
